# Update gh-pages to output generated at 456a3b4
# Increment "想去人数" (F column) counts by 1 for several events that
# appear on both the "展览" sheet and the aggregate "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F10").Value = 15
$wsExhibit.Range("F14").Value = 28
$wsExhibit.Range("F23").Value = 2308
$wsExhibit.Range("F32").Value = 1333

# Sheet "全部类型" (all types, aggregated)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F11").Value = 15
$wsAll.Range("F15").Value = 28
$wsAll.Range("F24").Value = 2308
$wsAll.Range("F33").Value = 1333
